$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 85, shifting existing rows 85-121 down to 86-122
$ws.Rows.Item(85).Insert()

# Populate the newly inserted row 85 with the new weekly price observation
$ws.Range("A85").Value = 7
$ws.Range("B85").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C85").Value = "Ñuble"
$ws.Range("D85").Value = 45202
$ws.Range("E85").Value = 16
$ws.Range("F85").Value = 100112044
$ws.Range("G85").Value = "Perejil"
$ws.Range("H85").Value = "Sin especificar"
$ws.Range("I85").Value = "Primera"
$ws.Range("J85").Value = 250
$ws.Range("K85").Value = 1500
$ws.Range("L85").Value = 1500
$ws.Range("M85").Value = 1500
$ws.Range("N85").Value = "`$/atado 0,5 a 1 kilo"
$ws.Range("O85").Value = "Región de Ñuble"
$ws.Range("P85").Value = 1500
$ws.Range("Q85").Value = 1
$ws.Range("R85").Value = "Hortaliza"
